$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(6, 8).Value = 175.14285
$ws.Cells.Item(6, 9).Value = 175.14285
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 525.4285500000001
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -413.4285500000001
$ws.Cells.Item(6, 14).ClearContents() | Out-Null
$ws.Cells.Item(11, 8).Value = 523.63635
$ws.Cells.Item(11, 9).Value = 523.63635
$ws.Cells.Item(11, 11).Value = 523.63635
$ws.Cells.Item(11, 13).Value = -383.63635
$ws.Cells.Item(20, 8).Value = 629.6667
$ws.Cells.Item(20, 9).Value = 629.6667
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 629.6667
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = -399.6667
$ws.Cells.Item(20, 14).ClearContents() | Out-Null
$ws.Cells.Item(33, 8).Value = 390
$ws.Cells.Item(33, 10).Value = 1073.3334
$ws.Cells.Item(33, 12).Value = 1073.3334
$ws.Cells.Item(33, 14).Value = -1531.3334
$ws.Cells.Item(35, 8).Value = 629.6667
$ws.Cells.Item(35, 9).Value = 629.6667
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 629.6667
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = -250.6667
$ws.Cells.Item(35, 14).ClearContents() | Out-Null
$ws.Cells.Item(86, 8).Value = 9750.5
$ws.Cells.Item(86, 9).Value = 8002
$ws.Cells.Item(86, 10).Value = 10333.333
$ws.Cells.Item(86, 11).Value = 8002
$ws.Cells.Item(86, 12).Value = 10333.333
$ws.Cells.Item(86, 13).Value = -6879
$ws.Cells.Item(86, 14).Value = -12579.333
$ws.Cells.Item(89, 8).Value = 9750.5
$ws.Cells.Item(89, 9).Value = 8002
$ws.Cells.Item(89, 10).Value = 10333.333
$ws.Cells.Item(89, 11).Value = 40010
$ws.Cells.Item(89, 12).Value = 51666.665
$ws.Cells.Item(89, 13).Value = -34394
$ws.Cells.Item(89, 14).Value = -62898.665

# Sheet ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 5560.5557
$ws.Cells.Item(32, 9).Value = 4620.8076
$ws.Cells.Item(32, 10).Value = 29994
$ws.Cells.Item(32, 11).Value = 4620.8076
$ws.Cells.Item(32, 12).Value = 29994
$ws.Cells.Item(32, 13).Value = -4333.8076
$ws.Cells.Item(32, 14).Value = -30568
$ws.Cells.Item(132, 8).Value = 3144.7144
$ws.Cells.Item(132, 10).Value = 4002.8
$ws.Cells.Item(132, 12).Value = 12008.4
$ws.Cells.Item(132, 14).Value = -17068.4

# Sheet BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(5, 8).Value = 1000
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 1000
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 1000
$ws.Cells.Item(5, 13).ClearContents() | Out-Null
$ws.Cells.Item(5, 14).Value = -1226
$ws.Cells.Item(105, 8).Value = 2685.25
$ws.Cells.Item(105, 9).Value = 2140.4285
$ws.Cells.Item(105, 11).Value = 2140.4285
$ws.Cells.Item(105, 13).Value = -393.4285

# Sheet CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(16, 8).Value = 964.1818
$ws.Cells.Item(16, 9).Value = 770.6667
$ws.Cells.Item(16, 11).Value = 770.6667
$ws.Cells.Item(16, 13).Value = -483.6667
$ws.Cells.Item(31, 8).Value = 3083.2632
$ws.Cells.Item(31, 9).Value = 2828.125
$ws.Cells.Item(31, 11).Value = 2828.125
$ws.Cells.Item(31, 13).Value = -2533.125
$ws.Cells.Item(34, 8).Value = 3083.2632
$ws.Cells.Item(34, 9).Value = 2828.125
$ws.Cells.Item(34, 11).Value = 2828.125
$ws.Cells.Item(34, 13).Value = -2626.125
$ws.Cells.Item(58, 8).Value = 8517.546
$ws.Cells.Item(58, 9).Value = 8854.333000000001
$ws.Cells.Item(58, 10).Value = 7002
$ws.Cells.Item(58, 11).Value = 8854.333000000001
$ws.Cells.Item(58, 12).Value = 7002
$ws.Cells.Item(58, 13).Value = -8651.333000000001
$ws.Cells.Item(58, 14).Value = -7408
$ws.Cells.Item(113, 8).Value = 964.1818
$ws.Cells.Item(113, 9).Value = 770.6667
$ws.Cells.Item(113, 11).Value = 770.6667
$ws.Cells.Item(113, 13).Value = 1399.3333
$ws.Cells.Item(132, 8).Value = 2547.762
$ws.Cells.Item(132, 9).Value = 1720.8572
$ws.Cells.Item(132, 10).Value = 4201.5713
$ws.Cells.Item(132, 11).Value = 5162.571599999999
$ws.Cells.Item(132, 12).Value = 12604.7139
$ws.Cells.Item(132, 13).Value = -2632.571599999999
$ws.Cells.Item(132, 14).Value = -17664.7139
$ws.Cells.Item(136, 8).Value = 8517.546
$ws.Cells.Item(136, 9).Value = 8854.333000000001
$ws.Cells.Item(136, 10).Value = 7002
$ws.Cells.Item(136, 11).Value = 26562.999
$ws.Cells.Item(136, 12).Value = 21006
$ws.Cells.Item(136, 13).Value = -24012.999
$ws.Cells.Item(136, 14).Value = -26106

# Sheet CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(4, 8).Value = 140.625
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 14).ClearContents() | Out-Null
$ws.Cells.Item(34, 8).Value = 3039.6
$ws.Cells.Item(34, 10).Value = 4966.3335
$ws.Cells.Item(34, 12).Value = 14899.0005
$ws.Cells.Item(34, 14).Value = -15067.0005
$ws.Cells.Item(47, 8).Value = 293.75
$ws.Cells.Item(47, 9).Value = 293.75
$ws.Cells.Item(47, 11).Value = 881.25
$ws.Cells.Item(47, 13).Value = -450.25
$ws.Cells.Item(55, 8).Value = 2593.1667
$ws.Cells.Item(55, 10).Value = 3098.4
$ws.Cells.Item(55, 12).Value = 9295.200000000001
$ws.Cells.Item(55, 14).Value = -9649.200000000001
$ws.Cells.Item(81, 8).Value = 3013
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 14).ClearContents() | Out-Null
$ws.Cells.Item(84, 8).Value = 3013
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 14).ClearContents() | Out-Null
$ws.Cells.Item(86, 8).Value = 1142.5
$ws.Cells.Item(86, 9).Value = 785.7143
$ws.Cells.Item(86, 10).Value = 1499.2858
$ws.Cells.Item(86, 11).Value = 2357.1429
$ws.Cells.Item(86, 12).Value = 4497.857400000001
$ws.Cells.Item(86, 13).Value = -1171.1429
$ws.Cells.Item(86, 14).Value = -6869.857400000001
$ws.Cells.Item(89, 8).Value = 1142.5
$ws.Cells.Item(89, 9).Value = 785.7143
$ws.Cells.Item(89, 10).Value = 1499.2858
$ws.Cells.Item(89, 11).Value = 7071.428699999999
$ws.Cells.Item(89, 12).Value = 13493.5722
$ws.Cells.Item(89, 13).Value = -1143.428699999999
$ws.Cells.Item(89, 14).Value = -25349.5722
$ws.Cells.Item(92, 8).Value = 894.05884
$ws.Cells.Item(92, 9).Value = 700
$ws.Cells.Item(92, 11).Value = 2100
$ws.Cells.Item(92, 13).Value = -852
$ws.Cells.Item(122, 8).Value = 1192.2142
$ws.Cells.Item(122, 10).Value = 1465.8889
$ws.Cells.Item(122, 12).Value = 13193.0001
$ws.Cells.Item(122, 14).Value = -18093.0001

# Sheet GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(10, 8).Value = 1756666.4
$ws.Cells.Item(10, 9).Value = 2509999.5
$ws.Cells.Item(10, 11).Value = 2509999.5
$ws.Cells.Item(10, 13).Value = -2509830.5
$ws.Cells.Item(14, 8).Value = 450974.9
$ws.Cells.Item(14, 9).Value = 501049.38
$ws.Cells.Item(14, 11).Value = 501049.38
$ws.Cells.Item(14, 13).Value = -500881.38

# Sheet LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(10, 8).Value = 167455.33
$ws.Cells.Item(10, 10).Value = 833
$ws.Cells.Item(10, 12).Value = 833
$ws.Cells.Item(10, 14).Value = -1113
$ws.Cells.Item(82, 8).Value = 3666.2666
$ws.Cells.Item(82, 9).Value = 1323.75
$ws.Cells.Item(82, 11).Value = 1323.75
$ws.Cells.Item(82, 13).Value = -962.75
$ws.Cells.Item(85, 8).Value = 3666.2666
$ws.Cells.Item(85, 9).Value = 1323.75
$ws.Cells.Item(85, 11).Value = 1323.75
$ws.Cells.Item(85, 13).Value = -75.75
$ws.Cells.Item(93, 8).Value = 792.5
$ws.Cells.Item(93, 9).Value = 792.5
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 792.5
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 455.5
$ws.Cells.Item(93, 14).ClearContents() | Out-Null

# Sheet WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 14).ClearContents() | Out-Null
$ws.Cells.Item(100, 8).Value = 9958824
$ws.Cells.Item(100, 9).Value = 34848884
$ws.Cells.Item(100, 11).Value = 69697768
$ws.Cells.Item(100, 13).Value = -69697227
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 14).ClearContents() | Out-Null
$ws.Cells.Item(122, 8).Value = 3958.6
$ws.Cells.Item(122, 9).Value = 3950.75
$ws.Cells.Item(122, 11).Value = 11852.25
$ws.Cells.Item(122, 13).Value = -9402.25
